# Apply the GitHub Actions "Updated cryptos list" data refresh
# to the crypto price/volume table on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D (Price) and E (Volume(1h)) columns hold text that can look numeric
# (e.g. "0.997", "56.148.19"). Force them to Text format before writing
# so Excel keeps them as strings instead of silently parsing them into
# numbers/dates, then restore the default "Normal" style afterwards so
# no stray per-cell formatting is left behind.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

# --- Row 47-49: coins rotated (VeChain / Bittensor / RenderToken) ---
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'

# --- Price / Volume(1h) refresh for every row ---
$ws.Range("D2").Value = '56.148.19'
$ws.Range("E2").Value = '  +0.50%  '
$ws.Range("D3").Value = '2.390.52'
$ws.Range("E3").Value = '  -4.43%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '479.51'
$ws.Range("E5").Value = '  -1.43%  '
$ws.Range("D6").Value = '147.99'
$ws.Range("E6").Value = '  +2.01%  '
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("E8").Value = '  -2.64%  '
$ws.Range("D9").Value = '2.391.69'
$ws.Range("E9").Value = '  -5.16%  '
$ws.Range("D10").Value = '0.0978'
$ws.Range("E10").Value = '  -0.04%  '
$ws.Range("D11").Value = '5.51'
$ws.Range("E11").Value = '  -2.40%  '
$ws.Range("D12").Value = '0.325'
$ws.Range("E12").Value = '  -2.46%  '
$ws.Range("E13").Value = '  +0.92%  '
$ws.Range("D14").Value = '2.806.17'
$ws.Range("E14").Value = '  -4.53%  '
$ws.Range("D15").Value = '56.442.52'
$ws.Range("E15").Value = '  +0.95%  '
$ws.Range("D16").Value = '20.39'
$ws.Range("E16").Value = '  -3.46%  '
$ws.Range("D17").Value = '0.0000132'
$ws.Range("E17").Value = '  -2.80%  '
$ws.Range("D18").Value = '2.395.49'
$ws.Range("E18").Value = '  -4.79%  '
$ws.Range("E19").Value = '  +1.65%  '
$ws.Range("D20").Value = '313.97'
$ws.Range("E20").Value = '  -2.04%  '
$ws.Range("D21").Value = '9.75'
$ws.Range("E21").Value = '  -5.05%  '
$ws.Range("D22").Value = '0.997'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").Value = '5.67'
$ws.Range("E23").Value = '  -2.60%  '
$ws.Range("D24").Value = '56.87'
$ws.Range("E24").Value = '  -2.78%  '
$ws.Range("E25").Value = '  +0.15%  '
$ws.Range("D26").Value = '0.395'
$ws.Range("E26").Value = '  -3.97%  '
$ws.Range("D27").Value = '0.158'
$ws.Range("E27").Value = '  -5.54%  '
$ws.Range("D28").Value = '2.499.52'
$ws.Range("E28").Value = '  -4.55%  '
$ws.Range("D29").Value = '7.31'
$ws.Range("E29").Value = '  -2.64%  '
$ws.Range("D30").Value = '0.0₃0771'
$ws.Range("E30").Value = '  -2.10%  '
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("D32").Value = '149.11'
$ws.Range("E32").Value = '  +0.07%  '
$ws.Range("D33").Value = '17.94'
$ws.Range("E33").Value = '  -2.40%  '
$ws.Range("E34").Value = '  -0.50%  '
$ws.Range("D35").Value = '4.97'
$ws.Range("E35").Value = '  -5.11%  '
$ws.Range("D36").Value = '1.11'
$ws.Range("E36").Value = '  -2.85%  '
$ws.Range("D37").Value = '0.851'
$ws.Range("D38").Value = '3.59'
$ws.Range("E38").Value = '  -2.87%  '
$ws.Range("D39").Value = '33.58'
$ws.Range("E39").Value = '  -1.97%  '
$ws.Range("D40").Value = '1.35'
$ws.Range("E40").Value = '  +2.41%  '
$ws.Range("E41").Value = '  +0.29%  '
$ws.Range("D42").Value = '0.0545'
$ws.Range("E42").Value = '  -2.14%  '
$ws.Range("D43").Value = '3.39'
$ws.Range("E43").Value = '  -4.70%  '
$ws.Range("D44").Value = '0.588'
$ws.Range("E44").Value = '  -4.24%  '
$ws.Range("D45").Value = '0.0946'
$ws.Range("E45").Value = '  +4.00%  '
$ws.Range("E46").Value = '  +0.24%  '
$ws.Range("D47").Value = '0.0224'
$ws.Range("E47").Value = '  -1.09%  '
$ws.Range("D48").Value = '252.70'
$ws.Range("E48").Value = '  -4.98%  '
$ws.Range("D49").Value = '4.56'
$ws.Range("E49").Value = '  -5.50%  '
$ws.Range("D50").Value = '17.05'
$ws.Range("E50").Value = '  -3.24%  '
$ws.Range("D51").Value = '1.782.81'
$ws.Range("E51").Value = '  -8.67%  '

# Restore default styling on the data range (undo the temporary Text format).
$dataRange.Style = "Normal"

